$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell values (translated from Portuguese to English)
$ws.Range("B1").Value = "AVERAGE_EFFORT_PER_MERGE_WO_FT"
$ws.Range("C1").Value = "AVERAGE_EFFORT_PER_MERGE_WITH_FT"

# Update the selected range to match the new selection B1:C1
$ws.Range("B1:C1").Select()
